{"js": "// Phase 3 requirements doc update:\n//  - \"Stage 3\" / \"Project 2\" -> \"Phase 3\" / \"Phase 2\" in the intro sentence\n//  - \"The Reducers\" -> \"The reducers\" (lowercase the r)\n//  - \"Although clearly, the mapper process ...\" -> \"The mapper process ...\"\n//  - \"...do either mapper or reducer\" -> \"...do both mapper and reducer\"\n//  - \"You may use, Boost unit testing\" -> \"You may use Boost unit testing\" (drop stray comma)\n//  - Remove the \"Code formatting & comments should follow Google style guide: ...\" bullet entirely\n\nconst body = context.document.body;\n\n// Locate the four sentences that get literal text substitutions.\nconst introRange = body.search(\n  \"In Stage 3, you will be extending your Project 2 solution:\",\n  { matchCase: true }\n);\nconst reducersRange = body.search(\n  \"The Reducers can be created once all of the mappers are complete.\",\n  { matchCase: true }\n);\nconst althoughRange = body.search(\n  \"Although clearly, the mapper process will need to run the mapper function from\",\n  { matchCase: true }\n);\nconst eitherRange = body.search(\n  \"You may either have a single executable that can do either mapper or reducer\",\n  { matchCase: true }\n);\nconst commaRange = body.search(\"You may use, Boost\", { matchCase: true });\n\nintroRange.load(\"items\");\nreducersRange.load(\"items\");\nalthoughRange.load(\"items\");\neitherRange.load(\"items\");\ncommaRange.load(\"items\");\nawait context.sync();\n\nif (introRange.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the intro sentence, found \" + introRange.items.length);\n}\nif (reducersRange.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the reducers sentence, found \" + reducersRange.items.length);\n}\nif (althoughRange.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the 'Although clearly' sentence, found \" + althoughRange.items.length);\n}\nif (eitherRange.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the 'either mapper' sentence, found \" + eitherRange.items.length);\n}\nif (commaRange.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'You may use, Boost', found \" + commaRange.items.length);\n}\n\nintroRange.items[0].insertText(\n  \"In Phase 3, you will be extending your Phase 2 solution:\",\n  Word.InsertLocation.replace\n);\nreducersRange.items[0].insertText(\n  \"The reducers can be created once all of the mappers are complete.\",\n  Word.InsertLocation.replace\n);\nalthoughRange.items[0].insertText(\n  \"The mapper process will need to run the mapper function from\",\n  Word.InsertLocation.replace\n);\neitherRange.items[0].insertText(\n  \"You may either have a single executable that can do both mapper and reducer\",\n  Word.InsertLocation.replace\n);\ncommaRange.items[0].insertText(\"You may use Boost\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// Delete the whole \"Code formatting & comments should follow Google style guide\" bullet,\n// including its hyperlink run, leaving the following \"Errors, warnings, ...\" bullet intact.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet styleGuideParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Code formatting\") !== -1) {\n    styleGuideParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!styleGuideParagraph) {\n  throw new Error(\"Could not find the 'Code formatting & comments...' paragraph to delete\");\n}\nstyleGuideParagraph.delete();\nawait context.sync();\n", "ps1": "# Phase 3 requirements doc update:\n#  - \"Stage 3\" / \"Project 2\" -> \"Phase 3\" / \"Phase 2\" in the intro sentence\n#  - \"The Reducers\" -> \"The reducers\" (lowercase the r)\n#  - \"Although clearly, the mapper process ...\" -> \"The mapper process ...\"\n#  - \"...do either mapper or reducer\" -> \"...do both mapper and reducer\"\n#  - \"You may use, Boost unit testing\" -> \"You may use Boost unit testing\" (drop stray comma)\n#  - Remove the \"Code formatting & comments should follow Google style guide: ...\" bullet entirely\n\n$d = $word.ActiveDocument\n\nfunction Assert-SingleMatch($haystack, $needle) {\n    $matchCount = ([regex]::Matches($haystack, [regex]::Escape($needle))).Count\n    if ($matchCount -ne 1) {\n        throw \"Expected exactly 1 match for '$needle', found $matchCount\"\n    }\n}\n\nfunction Replace-OnceText($find, $replace) {\n    $full = $d.Content.Text\n    Assert-SingleMatch $full $find\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($find, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-OnceText \"In Stage 3, you will be extending your Project 2 solution:\" \"In Phase 3, you will be extending your Phase 2 solution:\"\nReplace-OnceText \"The Reducers can be created once all of the mappers are complete.\" \"The reducers can be created once all of the mappers are complete.\"\nReplace-OnceText \"Although clearly, the mapper process will need to run the mapper function from\" \"The mapper process will need to run the mapper function from\"\nReplace-OnceText \"You may either have a single executable that can do either mapper or reducer\" \"You may either have a single executable that can do both mapper and reducer\"\nReplace-OnceText \"You may use, Boost\" \"You may use Boost\"\n\n# Delete the entire \"Code formatting & comments...\" bullet paragraph (including its\n# hyperlink run), leaving the following \"Errors, warnings, ...\" bullet in place.\n$styleGuideParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Code formatting*\") {\n        $styleGuideParagraph = $p\n        break\n    }\n}\nif ($styleGuideParagraph -eq $null) {\n    throw \"Could not find the 'Code formatting & comments...' paragraph to delete\"\n}\n$styleGuideParagraph.Range.Delete()\n"}
